$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "36.371.64"
$ws.Range("E2").Value = "  +1.30%  "
$ws.Range("D3").Value = "1.941.30"
$ws.Range("E3").Value = "  -0.60%  "
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "242.87"
$ws.Range("E5").Value = "  +0.94%  "
$ws.Range("E6").Value = "  -1.08%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "58.16"
$ws.Range("E7").Value = "  -3.43%  "
$ws.Range("E8").Value = "  -0.07%  "
$ws.Range("E9").Value = "  -1.99%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "55.83"
$ws.Range("E10").Value = "  -0.81%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0830"
$ws.Range("E11").Value = "  +4.02%  "
$ws.Range("E12").Value = "  +0.95%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "21.50"
$ws.Range("E13").Value = "  -1.10%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.820"
$ws.Range("E14").Value = "  -3.08%  "
$ws.Range("B15").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C15").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D15").Value = "2.227.10"
$ws.Range("E15").Value = "  -0.77%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "13.62"
$ws.Range("E16").Value = "  -1.91%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "5.23"
$ws.Range("E17").Value = "  -2.49%  "
$ws.Range("D18").Value = "1.946.14"
$ws.Range("E18").Value = "  -0.54%  "
$ws.Range("D19").Value = "36.319.00"
$ws.Range("E19").Value = "  +1.24%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "69.52"
$ws.Range("E20").Value = "  -1.46%  "
$ws.Range("E21").Value = "  +1.16%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "228.74"
$ws.Range("E22").Value = "  -2.42%  "
$ws.Range("E23").Value = "  -2.45%  "
$ws.Range("E24").Value = "  +0.08%  "
$ws.Range("E25").Value = "  -3.57%  "
$ws.Range("E26").Value = "  +0.56%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "9.18"
$ws.Range("E27").Value = "  -5.51%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "161.44"
$ws.Range("E28").Value = "  +1.43%  "
$ws.Range("B29").Value = "EthereumClassic"
$ws.Range("C29").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "19.43"
$ws.Range("E29").Value = "  -1.11%  "
$ws.Range("B30").Value = "Kaspa"
$ws.Range("C30").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.129"
$ws.Range("E30").Value = "  +0.18%  "
$ws.Range("E31").Value = "  -1.39%  "
$ws.Range("E32").Value = "  +1.73%  "
$ws.Range("E33").Value = "  -3.45%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "0.0625"
$ws.Range("E34").Value = "  +2.01%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "4.25"
$ws.Range("E35").Value = "  -2.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "6.27"
$ws.Range("E36").Value = "  +0.62%  "
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "1.78"
$ws.Range("E38").Value = "  -1.98%  "
$ws.Range("E39").Value = "  -5.33%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.02"
$ws.Range("E40").Value = "  -0.11%  "
$ws.Range("E41").Value = "  -0.25%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.87"
$ws.Range("E42").Value = "  +0.86%  "
$ws.Range("E43").Value = "  -3.57%  "
$ws.Range("E44").Value = "  -0.43%  "
$ws.Range("E45").Value = "  +0.36%  "
$ws.Range("D46").Value = "1.349.95"
$ws.Range("E46").Value = "  +1.64%  "
$ws.Range("E47").Value = "  -4.48%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "87.60"
$ws.Range("E48").Value = "  -4.15%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "7.09"
$ws.Range("E49").Value = "  -4.61%  "
$ws.Range("E50").Value = "  +0.39%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "45.03"
$ws.Range("E51").Value = "  +3.42%  "

$ws.Range("D5").Style = "Normal"
$ws.Range("D7").Style = "Normal"
$ws.Range("D10").Style = "Normal"
$ws.Range("D11").Style = "Normal"
$ws.Range("D13").Style = "Normal"
$ws.Range("D14").Style = "Normal"
$ws.Range("D16").Style = "Normal"
$ws.Range("D17").Style = "Normal"
$ws.Range("D20").Style = "Normal"
$ws.Range("D22").Style = "Normal"
$ws.Range("D27").Style = "Normal"
$ws.Range("D28").Style = "Normal"
$ws.Range("D29").Style = "Normal"
$ws.Range("D30").Style = "Normal"
$ws.Range("D34").Style = "Normal"
$ws.Range("D35").Style = "Normal"
$ws.Range("D36").Style = "Normal"
$ws.Range("D38").Style = "Normal"
$ws.Range("D40").Style = "Normal"
$ws.Range("D42").Style = "Normal"
$ws.Range("D48").Style = "Normal"
$ws.Range("D49").Style = "Normal"
$ws.Range("D51").Style = "Normal"
